$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A holds numeric-looking event codes that must stay text (like the source data). ---
# Excel infers type from the literal assigned to .Value, so a numeric-looking string would
# become a Number. Route through a Text-formatted scratch cell + PasteSpecial(values) so the
# destination keeps its original (default) style while the content stays Text, then drop the
# scratch row so it leaves no trace.
$codes = @(
    '113', '115', '155', '210', '215', '220', '300', '330', '340', '342', '346', '348', '352', '355', '356', '357', '365', '455', '465', '535', '549', '560', '580', '591', '610', '620', '740', '750', '800', '813', '831', '850'
)
$scratch = $ws.Range("Z100")
$scratch.NumberFormat = "@"
for ($i = 0; $i -lt $codes.Length; $i++) {
    $scratch.Value = $codes[$i]
    $scratch.Copy() | Out-Null
    $ws.Cells.Item($i + 2, 1).PasteSpecial(-4163) | Out-Null
}
$ws.Rows.Item(100).Delete() | Out-Null

# --- Columns B-E: nom_eve / Esperado / Observado / valor p ---
# Row 2
$ws.Range("B2").Value = 'Desnutrici”n aguda en menores de 5 anos'
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 0.01

# Row 3
$ws.Range("B3").Value = 'Cancer en menores de 18 anos'
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0

# Row 4
$ws.Range("B4").Value = 'Cancer de la mama y cuello uterino'
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 0.07

# Row 5
$ws.Range("B5").Value = 'Dengue'
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 9
$ws.Range("E5").Value = 0

# Row 6
$ws.Range("B6").Value = 'Defectos congenitos'
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 8
$ws.Range("E6").Value = 0

# Row 7
$ws.Range("B7").Value = 'Dengue grave'
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1

# Row 8
$ws.Range("B8").Value = 'Agresiones por animales potencialmente transmisores de rabia'
$ws.Range("C8").Value = 40
$ws.Range("D8").Value = 37
$ws.Range("E8").Value = 0.06

# Row 9
$ws.Range("B9").Value = 'Hepatitis a'
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.37

# Row 10
$ws.Range("B10").Value = 'Hepatitis b, c y coinfeccion hepatitis b y delta'
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 4
$ws.Range("E10").Value = 0.02

# Row 11
$ws.Range("B11").Value = 'Enfermedades huerfanas - raras'
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 0.05

# Row 12
$ws.Range("B12").Value = 'Ira por virus nuevo'
$ws.Range("C12").Value = 6
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0

# Row 13
$ws.Range("B13").Value = 'Infeccion respiratoria aguda grave irag inusitada'
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1

# Row 14
$ws.Range("B14").Value = 'Infecciones de sitio quirurgico asociadas a procedimiento medico quirurgico'
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0.37

# Row 15
$ws.Range("B15").Value = 'Enfermedad transmitida por alimentos o agua (eta)'
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 1

# Row 16
$ws.Range("B16").Value = 'Intento de suicidio'
$ws.Range("C16").Value = 11
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 0.12

# Row 17
$ws.Range("B17").Value = 'Iad - infecciones asociadas a dispositivos - individual'
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0.15

# Row 18
$ws.Range("B18").Value = 'Intoxicaciones'
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = 0.13

# Row 19
$ws.Range("B19").Value = 'Leptospirosis'
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 0.37

# Row 20
$ws.Range("B20").Value = 'Malaria'
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("B21").Value = 'Meningitis bacteriana y enfermedad meningoc”cica'
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0

# Row 22
$ws.Range("B22").Value = 'Morbilidad materna extrema'
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 0.16

# Row 23
$ws.Range("B23").Value = 'Mortalidad perinatal y neonatal tardia'
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0.37

# Row 24
$ws.Range("B24").Value = 'Mortalidad por dengue'
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 1

# Row 25
$ws.Range("B25").Value = 'Vigilancia integrada de muertes en menores de cinco anos por infeccion respiratoria aguda - enfermedad diarreica aguda y/o desnutricion'
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 0

# Row 26
$ws.Range("B26").ClearContents()
$ws.Range("C26").ClearContents()
$ws.Range("D26").Value = 0
$ws.Range("E26").ClearContents()

# Row 27
$ws.Range("B27").Value = 'Parotiditis'
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0.37

# Row 28
$ws.Range("B28").Value = 'Sifilis congenita'
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 1

# Row 29
$ws.Range("B29").Value = 'Sifilis gestacional'
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 0.37

# Row 30
$ws.Range("B30").Value = 'Tos ferina'
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0

# Row 31
$ws.Range("B31").Value = 'Tuberculosis'
$ws.Range("C31").Value = 8
$ws.Range("D31").Value = 5
$ws.Range("E31").Value = 0.09

# Row 32
$ws.Range("B32").Value = 'Varicela individual'
$ws.Range("C32").Value = 7
$ws.Range("D32").Value = 4
$ws.Range("E32").Value = 0.09

# Row 33
$ws.Range("B33").Value = 'Vih/sida/mortalidad por sida'
$ws.Range("C33").Value = 6
$ws.Range("D33").Value = 11
$ws.Range("E33").Value = 0.02

# Row 34 ("850"/Vih-sida) no longer exists in the updated table; the refreshed "850" data now
# lives in row 33, so drop the old trailing row entirely.
$ws.Rows.Item(34).Delete() | Out-Null

Write-Output "edit applied"